$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, $text)
    $cell.Value2 = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell $ws.Cells.Item(2,4) "41.854.39"
Set-TextCell $ws.Cells.Item(2,5) "  +0.19%  "
Set-TextCell $ws.Cells.Item(3,4) "2.271.87"
Set-TextCell $ws.Cells.Item(3,5) "  +0.13%  "
Set-TextCell $ws.Cells.Item(4,5) "  +0.00%  "
Set-TextCell $ws.Cells.Item(5,4) "303.39"
Set-TextCell $ws.Cells.Item(5,5) "  +0.09%  "
Set-TextCell $ws.Cells.Item(6,4) "92.92"
Set-TextCell $ws.Cells.Item(6,5) "  -0.17%  "
Set-TextCell $ws.Cells.Item(7,4) "0.530"
Set-TextCell $ws.Cells.Item(7,5) "  +1.40%  "
Set-TextCell $ws.Cells.Item(8,5) "  -0.09%  "
Set-TextCell $ws.Cells.Item(9,5) "  -0.62%  "
Set-TextCell $ws.Cells.Item(10,4) "32.58"
Set-TextCell $ws.Cells.Item(10,5) "  +0.50%  "
Set-TextCell $ws.Cells.Item(11,4) "53.59"
Set-TextCell $ws.Cells.Item(11,5) "  -1.53%  "
Set-TextCell $ws.Cells.Item(12,5) "  -0.28%  "
Set-TextCell $ws.Cells.Item(13,5) "  -1.61%  "
Set-TextCell $ws.Cells.Item(14,4) "6.70"
Set-TextCell $ws.Cells.Item(14,5) "  +0.45%  "
Set-TextCell $ws.Cells.Item(15,4) "2.624.50"
Set-TextCell $ws.Cells.Item(15,5) "  +0.03%  "
Set-TextCell $ws.Cells.Item(16,5) "  +0.91%  "
Set-TextCell $ws.Cells.Item(17,4) "2.283.44"
Set-TextCell $ws.Cells.Item(17,5) "  +0.99%  "
Set-TextCell $ws.Cells.Item(18,4) "0.781"
Set-TextCell $ws.Cells.Item(18,5) "  +3.51%  "
Set-TextCell $ws.Cells.Item(19,4) "41.774.72"
Set-TextCell $ws.Cells.Item(19,5) "  +0.24%  "
Set-TextCell $ws.Cells.Item(20,4) "12.80"
Set-TextCell $ws.Cells.Item(20,5) "  +3.00%  "
Set-TextCell $ws.Cells.Item(21,4) "0.0₃0909"
Set-TextCell $ws.Cells.Item(21,5) "  -0.06%  "
Set-TextCell $ws.Cells.Item(22,5) "  +0.15%  "
Set-TextCell $ws.Cells.Item(23,4) "67.32"
Set-TextCell $ws.Cells.Item(23,5) "  +0.18%  "
Set-TextCell $ws.Cells.Item(24,4) "244.09"
Set-TextCell $ws.Cells.Item(24,5) "  +1.27%  "
Set-TextCell $ws.Cells.Item(25,4) "2.59"
Set-TextCell $ws.Cells.Item(25,5) "  +0.23%  "
Set-TextCell $ws.Cells.Item(26,4) "1.94"
Set-TextCell $ws.Cells.Item(26,5) "  +3.20%  "
Set-TextCell $ws.Cells.Item(27,5) "  -0.09%  "
Set-TextCell $ws.Cells.Item(28,4) "24.07"
Set-TextCell $ws.Cells.Item(28,5) "  +0.95%  "
Set-TextCell $ws.Cells.Item(29,4) "9.56"
Set-TextCell $ws.Cells.Item(29,5) "  -1.66%  "
Set-TextCell $ws.Cells.Item(30,5) "  -5.35%  "
Set-TextCell $ws.Cells.Item(31,4) "34.99"
Set-TextCell $ws.Cells.Item(31,5) "  +2.36%  "
Set-TextCell $ws.Cells.Item(32,4) "160.45"
Set-TextCell $ws.Cells.Item(32,5) "  +1.24%  "
Set-TextCell $ws.Cells.Item(33,4) "5.26"
Set-TextCell $ws.Cells.Item(33,5) "  +1.32%  "
Set-TextCell $ws.Cells.Item(34,5) "  -0.05%  "
Set-TextCell $ws.Cells.Item(35,5) "  +0.68%  "
Set-TextCell $ws.Cells.Item(36,5) "  -1.15%  "
Set-TextCell $ws.Cells.Item(37,5) "  +1.84%  "
Set-TextCell $ws.Cells.Item(38,4) "16.88"
Set-TextCell $ws.Cells.Item(38,5) "  +1.39%  "
Set-TextCell $ws.Cells.Item(39,5) "  +0.15%  "
Set-TextCell $ws.Cells.Item(40,5) "  +0.57%  "
Set-TextCell $ws.Cells.Item(41,5) "  +0.57%  "
Set-TextCell $ws.Cells.Item(42,4) "3.95"
Set-TextCell $ws.Cells.Item(42,5) "  -1.07%  "
Set-TextCell $ws.Cells.Item(46,4) "10.46"
Set-TextCell $ws.Cells.Item(46,5) "  +3.94%  "
Set-TextCell $ws.Cells.Item(47,5) "  +7.52%  "
Set-TextCell $ws.Cells.Item(48,4) "2.90"
Set-TextCell $ws.Cells.Item(48,5) "  -2.67%  "
Set-TextCell $ws.Cells.Item(49,4) "53.33"
Set-TextCell $ws.Cells.Item(49,5) "  +3.03%  "
Set-TextCell $ws.Cells.Item(50,4) "73.18"
Set-TextCell $ws.Cells.Item(50,5) "  +3.79%  "

Set-TextCell $ws.Cells.Item(43,2) "EnergySwap"
Set-TextCell $ws.Cells.Item(43,3) "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws.Cells.Item(43,4) "19.86"
Set-TextCell $ws.Cells.Item(43,5) "  -2.37%  "

Set-TextCell $ws.Cells.Item(44,2) "Maker"
Set-TextCell $ws.Cells.Item(44,3) "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell $ws.Cells.Item(44,4) "2.014.81"
Set-TextCell $ws.Cells.Item(44,5) "  -2.23%  "

Set-TextCell $ws.Cells.Item(51,2) "TrustWalletToken"
Set-TextCell $ws.Cells.Item(51,3) "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws.Cells.Item(51,4) "1.15"
Set-TextCell $ws.Cells.Item(51,5) "  +1.10%  "

